$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.402.50"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").Value = "1.824.50"
$ws.Range("E3").Value = "  +2.23%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'316.42"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").Value = "'0.4007"
$ws.Range("E8").Value = "  +6.52%  "

$ws.Range("D9").Value = "'0.07644"
$ws.Range("E9").Value = "  +2.65%  "

$ws.Range("D10").Value = "'41.90"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("E11").Value = "  +1.69%  "

$ws.Range("D12").Value = "'6.318"
$ws.Range("E12").Value = "  +3.75%  "

$ws.Range("D13").Value = "'7.639"
$ws.Range("E13").Value = "  +6.00%  "

$ws.Range("D14").Value = "'1.002"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "'20.90"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "1.828.09"
$ws.Range("E16").Value = "  +2.91%  "

$ws.Range("D17").Value = "'89.64"
$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "'0.06589"
$ws.Range("E19").Value = "  +2.20%  "

$ws.Range("E20").Value = "  +2.10%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'6.071"
$ws.Range("E22").Value = "  +3.08%  "

$ws.Range("D23").Value = "28.418.55"
$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("D24").Value = "'11.11"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("E25").Value = "  +6.93%  "

$ws.Range("D26").Value = "'2.454"
$ws.Range("E26").Value = "  +7.77%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'157.19"
$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.65"
$ws.Range("E28").Value = "  +2.04%  "

$ws.Range("D29").Value = "2.038.41"
$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("D30").Value = "'123.78"
$ws.Range("E30").Value = "  +3.09%  "

$ws.Range("D31").Value = "'0.1117"
$ws.Range("E31").Value = "  +6.38%  "

$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("D33").Value = "'5.664"
$ws.Range("E33").Value = "  +2.35%  "

$ws.Range("D34").Value = "'0.07406"
$ws.Range("E34").Value = "  +15.32%  "

$ws.Range("D35").Value = "'3.646"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").Value = "'0.2235"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("D37").Value = "'0.02339"
$ws.Range("E37").Value = "  +2.47%  "

$ws.Range("D38").Value = "'5.229"
$ws.Range("E38").Value = "  +4.41%  "

$ws.Range("E39").Value = "  +4.72%  "

$ws.Range("D40").Value = "'0.6262"
$ws.Range("E40").Value = "  +1.93%  "

$ws.Range("E41").Value = "  +1.99%  "

$ws.Range("D42").Value = "'1.178"
$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "'1.394"
$ws.Range("E44").Value = "  -3.48%  "

$ws.Range("D45").Value = "'13.44"
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("D47").Value = "'0.5839"
$ws.Range("E47").Value = "  +1.55%  "

$ws.Range("D48").Value = "'124.80"
$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("D49").Value = "'2.003"
$ws.Range("E49").Value = "  +4.19%  "

$ws.Range("D50").Value = "'1.202"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").Value = "'0.06893"
$ws.Range("E51").Value = "  +1.45%  "
